# Adds the new "Introduction / Scene Vision / Workflow / Research / ... /
# Materials" heading paragraphs (plus the "Plane manipulation..." note and
# trailing blank paragraphs) to the end of the document body, after the
# existing "Materials" entry of the Index list.

$d = $word.ActiveDocument

function Insert-HeadingParagraph {
    param($doc, [string]$HeadingText)

    $lastPara = $doc.Paragraphs.Last
    $paraRange = $lastPara.Range
    # Collapsed point immediately before the paragraph mark of the current
    # last paragraph -- inserting here appends a brand new, independent
    # paragraph after it without disturbing the existing paragraph's runs.
    $insertPoint = $doc.Range($paraRange.End - 1, $paraRange.End - 1)

    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>$HeadingText</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $insertPoint.InsertXML($xml)
}

function Insert-PlainParagraph {
    param($doc, [string]$Text)

    $lastPara = $doc.Paragraphs.Last
    $paraRange = $lastPara.Range
    $insertPoint = $doc.Range($paraRange.End - 1, $paraRange.End - 1)

    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>$Text</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $insertPoint.InsertXML($xml)
}

function Insert-EmptyHeadingParagraph {
    param($doc)

    $lastPara = $doc.Paragraphs.Last
    $paraRange = $lastPara.Range
    $insertPoint = $doc.Range($paraRange.End - 1, $paraRange.End - 1)

    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $insertPoint.InsertXML($xml)
}

Insert-HeadingParagraph $d "Introduction"
Insert-HeadingParagraph $d "Scene Vision"
Insert-HeadingParagraph $d "Workflow"
Insert-HeadingParagraph $d "Research"
Insert-PlainParagraph   $d "Plane manipulation for crafting armour in blender"
Insert-HeadingParagraph $d "Scene Design"
Insert-HeadingParagraph $d "Asset List"
Insert-HeadingParagraph $d "Modelling"
Insert-HeadingParagraph $d "Materials"
Insert-EmptyHeadingParagraph $d
Insert-EmptyHeadingParagraph $d
Insert-EmptyHeadingParagraph $d
Insert-EmptyHeadingParagraph $d

Write-Output "Done. ParaCount=$($d.Paragraphs.Count)"
